$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new rows before row 726 (shifts existing 726:751 down to 729:754)
$ws.Range("A726:A728").EntireRow.Insert()

# Shared/fixed column values for this product block
$mercadoId   = 11
$mercado     = "Vega Monumental Concepción"
$region      = "Bíobío"
$codreg      = 8
$tipo        = "Fruta"
$productoId  = 100108
$producto    = "Tropicales y subtropicales"
$categoriaId = 100108006
$categoria   = "Plátano"
$variedad    = "Sin especificar"
$unidad      = "$/caja 20 kilos"
$origen      = "Ecuador"
$kgUnidad    = 20

# New week of data (fecha serial 45041) for the three calidades
$newRows = @(
    @{ Row=726; Calidad="Maduro";         Volumen=100; Min=16000; Max=16000; Prom=16000; KgPrecio=800 },
    @{ Row=727; Calidad="Pintón";         Volumen=300; Min=17000; Max=17000; Prom=17000; KgPrecio=850 },
    @{ Row=728; Calidad="Primera Pintón"; Volumen=300; Min=18000; Max=18000; Prom=18000; KgPrecio=900 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = 45041
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.KgPrecio
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
